# Weekly data refresh: insert a new week's record ahead of the existing
# "Poroto granado" rows (Vega Modelo de Temuco), pushing the old rows down
# by one and back-filling the top with the newest observation
# (fecha 44588 / 2022-01-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 53..63 down to 54..64, inserting a blank row 53.
$ws.Rows(53).Insert()

# Populate the new row 53 with the latest week's figures.
$ws.Cells.Item(53, 1).Value()  = 10
$ws.Cells.Item(53, 2).Value()  = "Vega Modelo de Temuco"
$ws.Cells.Item(53, 3).Value()  = "La Araucanía"
$ws.Cells.Item(53, 4).Value()  = 44588
$ws.Cells.Item(53, 5).Value()  = 9
$ws.Cells.Item(53, 6).Value()  = 100112030
$ws.Cells.Item(53, 7).Value()  = "Poroto granado"
$ws.Cells.Item(53, 8).Value()  = "Sin especificar"
$ws.Cells.Item(53, 9).Value()  = "Primera"
$ws.Cells.Item(53, 10).Value() = 185
$ws.Cells.Item(53, 11).Value() = 28000
$ws.Cells.Item(53, 12).Value() = 28000
$ws.Cells.Item(53, 13).Value() = 28000
$ws.Cells.Item(53, 14).Value() = "$/saco 25 kilos"
$ws.Cells.Item(53, 15).Value() = "Región de La Araucanía"
$ws.Cells.Item(53, 16).Value() = 1120
$ws.Cells.Item(53, 17).Value() = 25
$ws.Cells.Item(53, 18).Value() = "Hortaliza"
